{"js": "// Replace each \"old\u00d7digit=\" expression with its new value throughout the document body.\n// The mapping below was derived from the diff: each left-hand side is a unique\n// multiplication expression appearing exactly once in the document.\nconst replacements = [\n  [\"669\u00d79=\", \"744\u00d75=\"],\n  [\"509\u00d72=\", \"986\u00d77=\"],\n  [\"929\u00d72=\", \"876\u00d74=\"],\n  [\"467\u00d78=\", \"899\u00d74=\"],\n  [\"265\u00d78=\", \"570\u00d78=\"],\n  [\"254\u00d79=\", \"382\u00d77=\"],\n  [\"176\u00d73=\", \"221\u00d74=\"],\n  [\"992\u00d72=\", \"503\u00d75=\"],\n  [\"605\u00d79=\", \"889\u00d78=\"],\n  [\"841\u00d78=\", \"709\u00d76=\"],\n  [\"564\u00d76=\", \"269\u00d77=\"],\n  [\"985\u00d73=\", \"138\u00d73=\"],\n  [\"322\u00d76=\", \"529\u00d73=\"],\n  [\"709\u00d74=\", \"663\u00d78=\"],\n  [\"491\u00d75=\", \"277\u00d79=\"],\n  [\"956\u00d78=\", \"367\u00d72=\"],\n  [\"474\u00d75=\", \"950\u00d75=\"],\n  [\"890\u00d79=\", \"451\u00d73=\"],\n  [\"225\u00d72=\", \"428\u00d73=\"],\n  [\"724\u00d76=\", \"399\u00d79=\"],\n  [\"672\u00d78=\", \"897\u00d73=\"],\n  [\"968\u00d74=\", \"386\u00d78=\"],\n  [\"507\u00d73=\", \"514\u00d74=\"],\n  [\"843\u00d79=\", \"198\u00d73=\"],\n  [\"981\u00d76=\", \"209\u00d79=\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Each left-hand multiplication expression appears exactly once in the\n# document body, so a whole-document Find/Replace (wdReplaceAll) for each\n# pair reproduces the diff exactly.\n$replacements = @(\n    @(\"669\u00d79=\", \"744\u00d75=\"),\n    @(\"509\u00d72=\", \"986\u00d77=\"),\n    @(\"929\u00d72=\", \"876\u00d74=\"),\n    @(\"467\u00d78=\", \"899\u00d74=\"),\n    @(\"265\u00d78=\", \"570\u00d78=\"),\n    @(\"254\u00d79=\", \"382\u00d77=\"),\n    @(\"176\u00d73=\", \"221\u00d74=\"),\n    @(\"992\u00d72=\", \"503\u00d75=\"),\n    @(\"605\u00d79=\", \"889\u00d78=\"),\n    @(\"841\u00d78=\", \"709\u00d76=\"),\n    @(\"564\u00d76=\", \"269\u00d77=\"),\n    @(\"985\u00d73=\", \"138\u00d73=\"),\n    @(\"322\u00d76=\", \"529\u00d73=\"),\n    @(\"709\u00d74=\", \"663\u00d78=\"),\n    @(\"491\u00d75=\", \"277\u00d79=\"),\n    @(\"956\u00d78=\", \"367\u00d72=\"),\n    @(\"474\u00d75=\", \"950\u00d75=\"),\n    @(\"890\u00d79=\", \"451\u00d73=\"),\n    @(\"225\u00d72=\", \"428\u00d73=\"),\n    @(\"724\u00d76=\", \"399\u00d79=\"),\n    @(\"672\u00d78=\", \"897\u00d73=\"),\n    @(\"968\u00d74=\", \"386\u00d78=\"),\n    @(\"507\u00d73=\", \"514\u00d74=\"),\n    @(\"843\u00d79=\", \"198\u00d73=\"),\n    @(\"981\u00d76=\", \"209\u00d79=\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    # 1 = wdFindContinue, 2 = wdReplaceAll\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
